# Add a new "Save" column (H) to the s_vals sheet, mirroring the
# formatting of the existing header cells and appending a 0 value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (font, borders, alignment) from the last header cell (G1)
# onto the new header cell (H1), then set its text.
$ws.Cells.Item(1, 7).Copy() | Out-Null
$ws.Cells.Item(1, 8).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Cells.Item(1, 8).Value = "Save"

# New data value for the Save column.
$ws.Cells.Item(2, 8).Value = 0
